# Updated IPS AIP hipo turnover
# Applies revised turnover % figures across several location sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Grafton Wisconsin
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Grafton Wisconsin")
$ws.Range("E2").Value = 0.0811
$ws.Range("E3").Value = 0.0811
$ws.Range("E4").Value = 0.0811
$ws.Range("H4").Value = 0.0312
$ws.Range("I4").Value = 0.0082
$ws.Range("J4").Value = 0.0396
$ws.Range("K4").Value = 0.0164
$ws.Range("L4").Value = 0.0252
$ws.Range("M4").Value = 0
$ws.Range("N4").Value = 0.0416
$ws.Range("O4").Value = 0.0135166666666667
$ws.Range("P4").Value = 0.0135166666666667
$ws.Range("Q4").Value = 0.0135166666666667
$ws.Range("R4").Value = 0.04055
$ws.Range("S4").Value = 0.0135166666666667
$ws.Range("T4").Value = 0.0135166666666667
$ws.Range("U4").Value = 0.0135166666666667
$ws.Range("V4").Value = 0.04055
$ws.Range("W4").Value = 0.1622
$ws.Range("E5").Value = 0.857142857142857
$ws.Range("E6").Value = 0.857142857142857
$ws.Range("E7").Value = 0.857142857142857
$ws.Range("M7").Value = 1
$ws.Range("N7").Value = 1
$ws.Range("O7").Value = 0.857142857142857
$ws.Range("P7").Value = 0.857142857142857
$ws.Range("Q7").Value = 0.857142857142857
$ws.Range("R7").Value = 0.857142857142857
$ws.Range("S7").Value = 0.857142857142857
$ws.Range("T7").Value = 0.857142857142857
$ws.Range("U7").Value = 0.857142857142857
$ws.Range("V7").Value = 0.857142857142857
$ws.Range("W7").Value = 0.857142857142857

# ---------------------------------------------------------------------------
# Guadalajara Mexico
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Guadalajara Mexico")
$ws.Range("E2").Value = 0.7692
$ws.Range("E3").Value = 0.7692
$ws.Range("E4").Value = 0.7692
$ws.Range("M4").Value = 0
$ws.Range("N4").Value = 0
$ws.Range("O4").Value = 0.1282
$ws.Range("P4").Value = 0.1282
$ws.Range("Q4").Value = 0.1282
$ws.Range("R4").Value = 0.3846
$ws.Range("S4").Value = 0.1282
$ws.Range("T4").Value = 0.1282
$ws.Range("U4").Value = 0.1282
$ws.Range("V4").Value = 0.3846
$ws.Range("W4").Value = 1.5384

# ---------------------------------------------------------------------------
# Hyderabad India
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Hyderabad India")
$ws.Range("E2").Value = 0.0309
$ws.Range("E3").Value = 0.0309
$ws.Range("E4").Value = 0.0309
$ws.Range("I4").Value = 0.0156
$ws.Range("J4").Value = 0.0156
$ws.Range("L4").Value = 0.0156
$ws.Range("M4").Value = 0
$ws.Range("N4").Value = 0.0154
$ws.Range("O4").Value = 0.00515
$ws.Range("P4").Value = 0.00515
$ws.Range("Q4").Value = 0.00515
$ws.Range("R4").Value = 0.01545
$ws.Range("S4").Value = 0.00515
$ws.Range("T4").Value = 0.00515
$ws.Range("U4").Value = 0.00515
$ws.Range("V4").Value = 0.01545
$ws.Range("W4").Value = 0.0618

# ---------------------------------------------------------------------------
# Khed Taluka India: clear E2/E3 values (keep formatting) and remove row 4
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Khed Taluka India")
$ws.Range("E2").ClearContents()
$ws.Range("E3").ClearContents()
$ws.Rows.Item(4).Delete()

# ---------------------------------------------------------------------------
# Black River Falls Wisconsin
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Black River Falls Wisconsin")
$ws.Range("E2").Value = 0.6667
$ws.Range("E3").Value = 0.6667
$ws.Range("E4").Value = 0.6667
$ws.Range("M4").Value = 0
$ws.Range("N4").Value = 1
$ws.Range("O4").Value = 0.111116666666667
$ws.Range("P4").Value = 0.111116666666667
$ws.Range("Q4").Value = 0.111116666666667
$ws.Range("R4").Value = 0.33335
$ws.Range("S4").Value = 0.111116666666667
$ws.Range("T4").Value = 0.111116666666667
$ws.Range("U4").Value = 0.111116666666667
$ws.Range("V4").Value = 0.33335
$ws.Range("W4").Value = 1.3334

# ---------------------------------------------------------------------------
# South Beloit Gardner St Illino
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("South Beloit Gardner St Illino")
$ws.Range("E2").Value = 0.303
$ws.Range("E3").Value = 0.303
$ws.Range("E4").Value = 0.303
$ws.Range("M4").Value = 0
$ws.Range("N4").Value = 0
$ws.Range("O4").Value = 0.0505
$ws.Range("P4").Value = 0.0505
$ws.Range("Q4").Value = 0.0505
$ws.Range("R4").Value = 0.1515
$ws.Range("S4").Value = 0.0505
$ws.Range("T4").Value = 0.0505
$ws.Range("U4").Value = 0.0505
$ws.Range("V4").Value = 0.1515
$ws.Range("W4").Value = 0.606

# ---------------------------------------------------------------------------
# Chicago Lasalle Illinois
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Chicago Lasalle Illinois")
$ws.Range("E2").Value = 0.7692
$ws.Range("E3").Value = 0.7692
$ws.Range("E4").Value = 0.7692
$ws.Range("M4").Value = 0
$ws.Range("N4").Value = 0
$ws.Range("O4").Value = 0.1282
$ws.Range("P4").Value = 0.1282
$ws.Range("Q4").Value = 0.1282
$ws.Range("R4").Value = 0.3846
$ws.Range("S4").Value = 0.1282
$ws.Range("T4").Value = 0.1282
$ws.Range("U4").Value = 0.1282
$ws.Range("V4").Value = 0.3846
$ws.Range("W4").Value = 1.5384
